$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsCESTR = $wb.Worksheets.Item("CESTR")

# --- "About" sheet ---

# Source name changes (B3): "The Sales Tax Clearinghouse" -> "Secretaria da Receita Federal do Brasil"
$wsAbout.Range("B3").Value = "Secretaria da Receita Federal do Brasil"

# Drop the old hyperlink entirely (was on B6, pointed at thestc.com FAQ page).
$wsAbout.Hyperlinks.Delete()

# Old supporting rows (B4 "n/a", B5 "FAQ", B7 question text) are no longer used.
$wsAbout.Range("B4").Clear()
$wsAbout.Range("B5").Clear()
$wsAbout.Range("B6").Clear()
$wsAbout.Range("B7").Clear()

# New hyperlink cell lives at B4, displaying the Brazilian tax-authority URL.
$wsAbout.Range("B4").Value = "http://receita.economia.gov.br/"
$wsAbout.Hyperlinks.Add($wsAbout.Range("B4"), "https://thestc.com/FAQ.stm", "", "", "http://receita.economia.gov.br/") | Out-Null
$wsAbout.Range("B4").Style = "Hyperlink"

# "Note:" block moves up from row 9 to row 6, with new note text on row 7 (col A).
$wsAbout.Range("A9").Clear()
$wsAbout.Range("A10").Clear()

$wsAbout.Range("A6").Value = "Note:"
$wsAbout.Range("A6").Font.Bold = $true

$wsAbout.Range("A7").Value = "There are at least 6 different sales taxes in Brazil: ICMS, ISS, IE, COFINS, IPI and PIS/PASEP. We are using the standard ICMS rate."

# New supporting value for the rate itself: label + number, referenced by the CESTR sheet.
$wsAbout.Range("B9").Value = "ICMS = "
$wsAbout.Range("C9").Value = 0.17

# Restore view state (selection parked on an empty cell, as in the authored file).
$wsAbout.Activate()
$wsAbout.Range("B22").Select()

# --- "CESTR" sheet ---

# The tax rate is now pulled from the About sheet instead of being a hard-coded constant.
$wsCESTR.Range("B2").Formula = "=About!C9"

$wsCESTR.Activate()
$wsCESTR.Range("B3").Select()

$wsAbout.Activate()
